$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row updates for rows whose Coin/Link stay the same but Price / Volume(1h)
# values changed (straightforward cell value updates).
# ---------------------------------------------------------------------------
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.634.58"
$ws.Range("E2").Value = "  +1.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.274.71"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.80"
$ws.Range("E5").Value = "  +8.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "266.89"
$ws.Range("E6").Value = "  +0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +5.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +3.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.62"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +1.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.46"
$ws.Range("E12").Value = "  +7.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  -1.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.72"
$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.898"
$ws.Range("E15").Value = "  +5.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.616.03"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.272.33"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.693.65"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.92"
$ws.Range("E20").Value = "  +1.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.29"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.41"
$ws.Range("E22").Value = "  -2.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.76"
$ws.Range("E23").Value = "  +1.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.91"
$ws.Range("E24").Value = "  +1.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.53"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.20"
$ws.Range("E26").Value = "  +8.02%  "

$ws.Range("E27").Value = "  +1.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.03"
$ws.Range("E28").Value = "  +4.46%  "

$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.44"
$ws.Range("E31").Value = "  +1.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.47"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0915"
$ws.Range("E33").Value = "  +1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  -1.28%  "

$ws.Range("E39").Value = "  +4.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").Value = "  -3.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.89"
$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.68"
$ws.Range("E42").Value = "  -3.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.239"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.77"
$ws.Range("E46").Value = "  -5.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "75.03"
$ws.Range("E47").Value = "  +42.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.59"
$ws.Range("E48").Value = "  +2.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.27"
$ws.Range("E49").Value = "  +2.08%  "

# ---------------------------------------------------------------------------
# Rows 36/37: VeChain and NEARProtocol swap places (NEARProtocol now ranks
# above VeChain), each with updated Price / Volume(1h) values.
# ---------------------------------------------------------------------------
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.23"
$ws.Range("E36").Value = "  +10.09%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0383"
$ws.Range("E37").Value = "  +8.69%  "

# ---------------------------------------------------------------------------
# Rows 50/51: TheSandbox and FraxShare swap places (FraxShare now ranks
# above TheSandbox), each with updated Price / Volume(1h) values.
# ---------------------------------------------------------------------------
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.55"
$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.655"
$ws.Range("E51").Value = "  +16.65%  "
